$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6673.706
$ws.Range("J74").Value = 7556.1904
$ws.Range("L74").Value = 7556.1904
$ws.Range("N74").Value = -9428.190399999999

$ws.Range("H77").Value = 6673.706
$ws.Range("J77").Value = 7556.1904
$ws.Range("L77").Value = 37780.952
$ws.Range("N77").Value = -47140.952

$ws.Range("H97").Value = 2871.8572
$ws.Range("J97").Value = 3150.5
$ws.Range("L97").Value = 9451.5
$ws.Range("N97").Value = -10443.5

$ws.Range("H100").Value = 8041.25
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H116").Value = 8501.286
$ws.Range("I116").Value = 12088.286
$ws.Range("K116").Value = 12088.286
$ws.Range("M116").Value = -8646.286

$ws.Range("H135").Value = 1368.5
$ws.Range("I135").Value = 1402.3
$ws.Range("K135").Value = 12620.7
$ws.Range("M135").Value = -10085.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1398.8
$ws.Range("I63").Value = 1248.75
$ws.Range("K63").Value = 1248.75
$ws.Range("M63").Value = -562.75

$ws.Range("H66").Value = 1398.8
$ws.Range("I66").Value = 1248.75
$ws.Range("K66").Value = 6243.75
$ws.Range("M66").Value = -2811.75

$ws.Range("H132").Value = 2446.28
$ws.Range("I132").Value = 2295.7673
$ws.Range("K132").Value = 6887.3019
$ws.Range("M132").Value = -4357.3019

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 21665.572
$ws.Range("J76").Value = 21665.572
$ws.Range("L76").Value = 21665.572
$ws.Range("N76").Value = -22295.572

$ws.Range("H79").Value = 21665.572
$ws.Range("J79").Value = 21665.572
$ws.Range("L79").Value = 21665.572
$ws.Range("N79").Value = -23849.572

$ws.Range("H86").Value = 13960.818
$ws.Range("I86").Value = 3639
$ws.Range("K86").Value = 3639
$ws.Range("M86").Value = -2516

$ws.Range("H89").Value = 13960.818
$ws.Range("I89").Value = 3639
$ws.Range("K89").Value = 18195
$ws.Range("M89").Value = -12579

$ws.Range("H134").Value = 6776.164
$ws.Range("I134").Value = 5801.34
$ws.Range("K134").Value = 17404.02
$ws.Range("M134").Value = -14869.02

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 34
$ws.Range("J7").Value = 85
$ws.Range("L7").Value = 85
$ws.Range("N7").Value = -311

$ws.Range("H31").Value = 1501.8793
$ws.Range("J31").Value = 2119.9644
$ws.Range("L31").Value = 2119.9644
$ws.Range("N31").Value = -2709.9644

$ws.Range("H34").Value = 1501.8793
$ws.Range("J34").Value = 2119.9644
$ws.Range("L34").Value = 2119.9644
$ws.Range("N34").Value = -2523.9644

$ws.Range("H58").Value = 6110.125
$ws.Range("I58").Value = 2447.5
$ws.Range("J58").Value = 7331
$ws.Range("K58").Value = 2447.5
$ws.Range("L58").Value = 7331
$ws.Range("M58").Value = -2244.5
$ws.Range("N58").Value = -7737

$ws.Range("H86").Value = 6045.0835
$ws.Range("I86").Value = 7012.75
$ws.Range("K86").Value = 7012.75
$ws.Range("M86").Value = -5889.75

$ws.Range("H89").Value = 6045.0835
$ws.Range("I89").Value = 7012.75
$ws.Range("K89").Value = 35063.75
$ws.Range("M89").Value = -29447.75

$ws.Range("H122").Value = 2433.5
$ws.Range("I122").Value = 2451.1177
$ws.Range("J122").Value = 2373.6
$ws.Range("K122").Value = 7353.353099999999
$ws.Range("L122").Value = 7120.799999999999
$ws.Range("M122").Value = -4903.353099999999
$ws.Range("N122").Value = -12020.8

$ws.Range("H132").Value = 2553.4194
$ws.Range("I132").Value = 2538.5334
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 7615.600199999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5085.600199999999
$ws.Range("N132").Value = -14060

$ws.Range("H136").Value = 6110.125
$ws.Range("I136").Value = 2447.5
$ws.Range("J136").Value = 7331
$ws.Range("K136").Value = 7342.5
$ws.Range("L136").Value = 21993
$ws.Range("M136").Value = -4792.5
$ws.Range("N136").Value = -27093

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3750.7778
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H113").Value = 840.64
$ws.Range("I113").Value = 755.75
$ws.Range("J113").Value = 919
$ws.Range("K113").Value = 2267.25
$ws.Range("L113").Value = 2757
$ws.Range("M113").Value = -97.25
$ws.Range("N113").Value = -7097

$ws.Range("H139").Value = 5845
$ws.Range("I139").Value = 2554.65
$ws.Range("K139").Value = 7663.950000000001
$ws.Range("M139").Value = -2523.950000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7456.1665
$ws.Range("I41").Value = 7456.1665
$ws.Range("K41").Value = 7456.1665
$ws.Range("M41").Value = -7101.1665

$ws.Range("H70").Value = 7527.4443
$ws.Range("I70").Value = 6968.5
$ws.Range("J70").Value = 11999
$ws.Range("K70").Value = 6968.5
$ws.Range("L70").Value = 11999
$ws.Range("M70").Value = -6698.5
$ws.Range("N70").Value = -12539

$ws.Range("H73").Value = 7527.4443
$ws.Range("I73").Value = 6968.5
$ws.Range("J73").Value = 11999
$ws.Range("K73").Value = 6968.5
$ws.Range("L73").Value = 11999
$ws.Range("M73").Value = -6032.5
$ws.Range("N73").Value = -13871

$ws.Range("H80").Value = 2333.3333
$ws.Range("J80").Value = 2750
$ws.Range("L80").Value = 2750
$ws.Range("N80").Value = -4746

$ws.Range("H83").Value = 2333.3333
$ws.Range("J83").Value = 2750
$ws.Range("L83").Value = 13750
$ws.Range("N83").Value = -23734

$ws.Range("H102").Value = 2553.682
$ws.Range("J102").Value = 2998.75
$ws.Range("L102").Value = 2998.75
$ws.Range("N102").Value = -6242.75

$ws.Range("H132").Value = 3731.34
$ws.Range("I132").Value = 3674.3901
$ws.Range("J132").Value = 3990.7778
$ws.Range("K132").Value = 11023.1703
$ws.Range("L132").Value = 11972.3334
$ws.Range("M132").Value = -8493.1703
$ws.Range("N132").Value = -17032.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2004
$ws.Range("J100").Value = 2004
$ws.Range("L100").Value = 2004
$ws.Range("N100").Value = -3086

$ws.Range("H130").Value = 89999
$ws.Range("J130").Value = 89999
$ws.Range("L130").Value = 89999
$ws.Range("N130").Value = -100039

$ws.Range("H132").Value = 2468.0488
$ws.Range("I132").Value = 1858.037
$ws.Range("J132").Value = 3644.5
$ws.Range("K132").Value = 5574.111
$ws.Range("L132").Value = 10933.5
$ws.Range("M132").Value = -3044.111
$ws.Range("N132").Value = -15993.5

$ws.Range("H136").Value = 2427
$ws.Range("I136").Value = 1464.5
$ws.Range("J136").Value = 2977
$ws.Range("K136").Value = 4393.5
$ws.Range("L136").Value = 8931
$ws.Range("M136").Value = -1843.5
$ws.Range("N136").Value = -14031

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1452.7858
$ws.Range("I126").Value = 1395.7142
$ws.Range("K126").Value = 4187.142599999999
$ws.Range("M126").Value = -1717.142599999999
